$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($target, $donor, $text) {
    $ws.Range($target).Value = "'" + $text
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($target).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

function Set-NumCell($target, $donor, $num) {
    $ws.Range($target).Value = $num
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($target).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# --- Header text updates (Volume/Number and week-of dates) ---
$ws.Range("A8").Value = "Volume 32   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/20/2025  Through  1/26/2025"

# --- Data table updates (rows 15-31) ---
$ws.Range("N15").Value = -50
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 10
$ws.Range("J16").Value = 10
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -58.333333333333
$ws.Range("N16").Value = -87.80487804878
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 13
$ws.Range("K17").Value = -23.076923076923
$ws.Range("L17").Value = -50
$ws.Range("M17").Value = 150
$ws.Range("N17").Value = -69.696969696969
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 118.181818181818
$ws.Range("I18").Value = 20
$ws.Range("J18").Value = 10
$ws.Range("K18").Value = 100
$ws.Range("L18").Value = 25
$ws.Range("M18").Value = 566.666666666667
$ws.Range("N18").Value = -63.636363636363
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -35.714285714285
$ws.Range("F19").Value = 46
$ws.Range("H19").Value = 4.545454545454
$ws.Range("I19").Value = 42
$ws.Range("J19").Value = 43
$ws.Range("K19").Value = -2.325581395348
$ws.Range("L19").Value = -10.63829787234
$ws.Range("M19").Value = 68
$ws.Range("N19").Value = -41.666666666666
Set-TextCell "C20" "C14" "0"
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 400
$ws.Range("L20").Value = 33.333333333333
$ws.Range("M20").Value = -20
$ws.Range("N20").Value = -92
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -9.090909090909
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = 19.277108433734
$ws.Range("I21").Value = 87
$ws.Range("J21").Value = 78
$ws.Range("K21").Value = 11.538461538461
$ws.Range("L21").Value = -6.451612903225
$ws.Range("M21").Value = 42.622950819672
$ws.Range("N21").Value = -70.408163265306
Set-TextCell "D22" "C14" "0"
Set-TextCell "E22" "C14" "***.*"
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 3
$ws.Range("K22").Value = 50
$ws.Range("L22").Value = -25
$ws.Range("M22").Value = -66.666666666666
Set-TextCell "C23" "C14" "0"
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -50
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = -50
$ws.Range("L23").Value = -33.333333333333
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -17.142857142857
$ws.Range("F24").Value = 112
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = -9.677419354838
$ws.Range("I24").Value = 100
$ws.Range("J24").Value = 104
$ws.Range("K24").Value = -3.846153846153
$ws.Range("L24").Value = -19.354838709677
$ws.Range("M24").Value = -20
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 35
$ws.Range("E25").Value = -37.142857142857
$ws.Range("F25").Value = 76
$ws.Range("G25").Value = 110
$ws.Range("H25").Value = -30.90909090909
$ws.Range("I25").Value = 70
$ws.Range("J25").Value = 94
$ws.Range("K25").Value = -25.531914893617
$ws.Range("L25").Value = -42.622950819672
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -20
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 20.689655172413
$ws.Range("I26").Value = 31
$ws.Range("J26").Value = 28
$ws.Range("K26").Value = 10.714285714285
$ws.Range("L26").Value = 6.896551724137
$ws.Range("M26").Value = 0
Set-TextCell "C27" "C14" "0"
Set-NumCell "C28" "C26" 4
Set-NumCell "D28" "D26" 1
Set-NumCell "E28" "E26" 300
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 6
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = 100
$ws.Range("L28").Value = 50
Set-TextCell "D31" "C14" "0"
Set-TextCell "E31" "C14" "***.*"
